# Update "reporte_diario" style sheet:
#   - Remove the trailing rows (24-31) so the data range shrinks to A1:E23
#   - Update the report date (column A) from 45822 to 45825 for all data rows
#   - Replace the diet names (column C) and quantities (column D) with the
#     new, filtered ("por CDS") data set for the new date

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, delete the now-unused rows 24-31 so the sheet shrinks to A1:E23.
$ws.Range("A24:E31").EntireRow.Delete() | Out-Null

# New dataset for rows 2-23 (Fecha, Servicio, Dieta, Cantidad, Valor Total)
$data = @(
    @(45825, "Desayuno", "Astringente", 3, 0),
    @(45825, "Desayuno", "Blanda", 19, 0),
    @(45825, "Desayuno", "Coronaria", 15, 0),
    @(45825, "Desayuno", "Hepatica", 1, 0),
    @(45825, "Desayuno", "Hipercalorica", 1, 0),
    @(45825, "Desayuno", "Hiperproteica", 2, 0),
    @(45825, "Desayuno", "Hipo Grasa", 6, 0),
    @(45825, "Desayuno", "Hipoglucida", 12, 0),
    @(45825, "Desayuno", "Hiposodica", 28, 0),
    @(45825, "Desayuno", "Liquida Clara", 4, 0),
    @(45825, "Desayuno", "Liquida Total", 7, 0),
    @(45825, "Desayuno", "Liquida Total 140 Cc", 1, 0),
    @(45825, "Desayuno", "Liquida Total Miel 140 Cc", 2, 0),
    @(45825, "Desayuno", "Liquida Total Nectar", 9, 0),
    @(45825, "Desayuno", "Liquida Total Nectar 140 Cc", 2, 0),
    @(45825, "Desayuno", "Liquida total Miel", 3, 0),
    @(45825, "Desayuno", "Normal", 57, 0),
    @(45825, "Desayuno", "Renal Dialisis", 4, 0),
    @(45825, "Desayuno", "Renal PRE Dialisis", 5, 0),
    @(45825, "Desayuno", "Semiblanda", 28, 0),
    @(45825, "Desayuno", "Semiblanda Pequena", 4, 0),
    @(45825, "Desayuno", "Todo Pure", 2, 0)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
